# "double digits lesson implemented"
# This script updates the Key/Value language table on sheet "en":
#  - Rewords several lesson1 attack strings (rows 35-42)
#  - Inserts a brand new "lesson4" block (double digit multiplication lesson)
#    as rows 46-55, pushing the bonus-blob / bonus-distribute / bonus-partial
#    rows down by 10 (they become rows 56-66)
#  - Updates the sheet view's selection to reflect the new bottom of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Reword the existing lesson1 "attack" explanation strings (rows 35-42)
#    Column A (the Key) is unchanged for these rows - only the Value (B)
#    text is updated. Row 36's Value cell additionally gets vertical-center
#    alignment (matching the styling used by other multi-line Value cells).
# ---------------------------------------------------------------------
$ws.Range("B35").Value = "Here we will be splitting the numbers up by multiples of 10’s."
$ws.Range("B36").Value = "Splitting up the numbers this way allows us to simply multiply single digit numbers, and just add the zeroes in the end."
$ws.Range("B36").VerticalAlignment = -4108
$ws.Range("B37").Value = "Now click on the boxed number to split it!"
$ws.Range("B38").Value = "Here we will compute the product of each sub areas by using the numpad to type in the digits."
$ws.Range("B39").Value = "Since we split them up nicely, we only ever have to worry about multiplying single digit numbers."
$ws.Range("B40").Value = "Just remember to add the zero at the end for the number with the denomination of 10!"
$ws.Range("B42").Value = "Finally, we add the numbers together to form the final product!"

# ---------------------------------------------------------------------
# 2. Insert 10 new blank rows before the old row 46 (bonusBlob_1), shifting
#    all the bonus* rows down so the new lesson4 content can be placed
#    right after the lesson1 content, keeping the table logically grouped.
# ---------------------------------------------------------------------
$ws.Rows("46:55").Insert()

# ---------------------------------------------------------------------
# 3. Populate the new lesson4 (double digits) rows.
#    Column A = Key, Column B = Value. Some Value cells use vertical-center
#    alignment (style index 3), matching the same convention used
#    elsewhere in the sheet for longer / multi-line dialogue strings.
# ---------------------------------------------------------------------
$lesson4 = @(
  @(46, "lesson4_intro_1",             "Here we are at the final stage, with only a handful of blobs left to banish!", $false),
  @(47, "lesson4_intro_2",             "This time around, we will be multiplying double digit numbers.", $true),
  @(48, "lesson4_area_1",              "Just as we split the area up horizontally, we can also split it vertically.", $true),
  @(49, "lesson4_connect_1",           "Now why don’t you give it a try?", $true),
  @(50, "lesson4_attack_distribute_1", "Just as you have done up to this point, simply click on the boxed numbers to split them.", $false),
  @(51, "lesson4_attack_eval_dd_1",    "This time around, you will be multiplying two double digit numbers.", $true),
  @(52, "lesson4_attack_eval_dd_2",    "Fortunately, they are both multiples of 10.", $true),
  @(53, "lesson4_attack_eval_dd_3",    "All you need to do is multiply the two non-zero digits, and put two zeroes at the end.", $true),
  @(54, "lesson4_end_1",               "Nice! Even with multiplying double digits, splitting them up allows you to solve them with ease!", $false),
  @(55, "lesson4_end_2",               "You will certainly have no problem dealing with the remaining blobs!", $false)
)

foreach ($entry in $lesson4) {
    $row = $entry[0]
    $key = $entry[1]
    $value = $entry[2]
    $center = $entry[3]

    $ws.Cells.Item($row, 1).Value = $key
    $bcell = $ws.Cells.Item($row, 2)
    $bcell.Value = $value
    if ($center) {
        $bcell.VerticalAlignment = -4108
    }
}

# ---------------------------------------------------------------------
# 4. Update the sheet view to match the new selection / scroll position.
# ---------------------------------------------------------------------
$ws.Range("A55").Select()
$excel.ActiveWindow.ScrollRow = 28
